# Applies the updated Betfair odds values to Sheet1.
# The workbook is already open; grab the active workbook/sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("F2").Value = 3.95
$ws.Range("H2").Value = 1.65
$ws.Range("K2").Value = 5
$ws.Range("P2").Value = 2.12

# Row 3
$ws.Range("Q3").Value = 2.9

# Row 4
$ws.Range("F4").Value = 1.59
$ws.Range("G4").Value = 1.73

# Row 5
$ws.Range("F5").Value = 1.8
$ws.Range("G5").Value = 1.96
$ws.Range("H5").Value = 4.7
$ws.Range("I5").Value = 6
$ws.Range("J5").Value = 3.45
$ws.Range("K5").Value = 4.1
$ws.Range("M5").Value = 1.08
$ws.Range("P5").Value = 1.81
$ws.Range("Q5").Value = 2
$ws.Range("T5").Value = 1.87
$ws.Range("U5").Value = 1.9
$ws.Range("AA5").Value = 150
$ws.Range("AB5").Value = 9.199999999999999
$ws.Range("AC5").Value = 9.6
$ws.Range("AF5").Value = 1000
$ws.Range("AG5").Value = 1000
$ws.Range("AI5").Value = 90
$ws.Range("AM5").Value = 150

# Row 7
$ws.Range("J7").Value = 8.199999999999999
$ws.Range("P7").Value = 2.8
$ws.Range("Q7").Value = 1.32

# Row 8
$ws.Range("F8").Value = 2.62
$ws.Range("G8").Value = 3.05
$ws.Range("H8").Value = 3.35
$ws.Range("I8").Value = 4.1
$ws.Range("K8").Value = 2.98

# Row 9
$ws.Range("F9").Value = 1.8
$ws.Range("G9").Value = 1.99
$ws.Range("H9").Value = 5.3
$ws.Range("J9").Value = 2.8
$ws.Range("K9").Value = 3.6
$ws.Range("N9").Value = 2.76

# Row 10
$ws.Range("F10").Value = 1.8
$ws.Range("I10").Value = 6.8
$ws.Range("K10").Value = 4.2
$ws.Range("P10").Value = 1.69
$ws.Range("Q10").Value = 2.18

# Row 14
$ws.Range("F14").Value = 1.69
$ws.Range("H14").Value = 5.7
